$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 29 (shifts existing rows 29..133 down to 30..134)
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly price record
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44592
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 100112030
$ws.Range("G29").Value = "Poroto granado"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 35
$ws.Range("K29").Value = 23000
$ws.Range("L29").Value = 23000
$ws.Range("M29").Value = 23000
$ws.Range("N29").Value = "$/malla 25 kilos"
$ws.Range("O29").Value = "Provincia de Talca"
$ws.Range("P29").Value = 920
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
